# Auto-generated edit script: updates Leve profit/cost figures across all sheets
# per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (61 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I8").Value = 62.75
$ws.Range("N8").Value = -3953
$ws.Range("J8").Value = 1225
$ws.Range("L8").Value = 3675
$ws.Range("H8").Value = 450.16666
$ws.Range("K8").Value = 188.25
$ws.Range("M8").Value = -49.25
$ws.Range("K80").Value = 1806
$ws.Range("I80").Value = 602
$ws.Range("L80").Value = 1627.5
$ws.Range("N80").Value = -3623.5
$ws.Range("J80").Value = 542.5
$ws.Range("M80").Value = -808
$ws.Range("H80").Value = 582.1667
$ws.Range("H82").Value = 1858.4
$ws.Range("M82").Value = -5169.200000000001
$ws.Range("K82").Value = 5575.200000000001
$ws.Range("I82").Value = 1858.4
$ws.Range("M83").Value = -426
$ws.Range("I83").Value = 602
$ws.Range("L83").Value = 4882.5
$ws.Range("N83").Value = -14866.5
$ws.Range("H83").Value = 582.1667
$ws.Range("J83").Value = 542.5
$ws.Range("K83").Value = 5418
$ws.Range("I85").Value = 1858.4
$ws.Range("M85").Value = -4171.200000000001
$ws.Range("H85").Value = 1858.4
$ws.Range("K85").Value = 5575.200000000001
$ws.Range("N92").Value = -133336706
$ws.Range("H92").Value = 66138256
$ws.Range("M92").Value = -5049780.5
$ws.Range("J92").Value = 133334210
$ws.Range("K92").Value = 5051028.5
$ws.Range("I92").Value = 5051028.5
$ws.Range("L92").Value = 133334210
$ws.Range("K96").Value = 2388
$ws.Range("I96").Value = 796
$ws.Range("H96").Value = 808.2
$ws.Range("M96").Value = -1015
$ws.Range("N96").Value = -5194.9999
$ws.Range("J96").Value = 816.3333
$ws.Range("L96").Value = 2448.9999
$ws.Range("J129").Value = 1195.0857
$ws.Range("M129").Value = 3206
$ws.Range("I129").Value = 598
$ws.Range("K129").Value = 1794
$ws.Range("N129").Value = -13585.2571
$ws.Range("L129").Value = 3585.2571
$ws.Range("H129").Value = 1120.45
$ws.Range("N131").Value = -19397.3688
$ws.Range("L131").Value = 9317.3688
$ws.Range("J131").Value = 3105.7896
$ws.Range("H131").Value = 2646.0417
$ws.Range("H137").Value = 1658
$ws.Range("K137").Value = 5344.125
$ws.Range("L137").Value = 3000
$ws.Range("I137").Value = 1781.375
$ws.Range("M137").Value = -2794.125
$ws.Range("N137").Value = -8100
$ws.Range("J137").Value = 1000

# --- Sheet: ARM (52 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N8").Value = -10288
$ws.Range("J8").Value = 10000
$ws.Range("L8").Value = 10000
$ws.Range("H8").Value = 10000
$ws.Range("N21").ClearContents()
$ws.Range("H21").Value = 5010
$ws.Range("J21").Value = 0
$ws.Range("M21").Value = -4636
$ws.Range("L21").Value = 0
$ws.Range("I21").Value = 5010
$ws.Range("K21").Value = 5010
$ws.Range("N27").ClearContents()
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("I32").Value = 5179.634
$ws.Range("M32").Value = -4892.634
$ws.Range("H32").Value = 7182.1025
$ws.Range("K32").Value = 5179.634
$ws.Range("N34").ClearContents()
$ws.Range("L34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L61").Value = 3987.6843
$ws.Range("J61").Value = 3987.6843
$ws.Range("N61").Value = -4411.6843
$ws.Range("I61").Value = 6562.3184
$ws.Range("H61").Value = 5369.1953
$ws.Range("M61").Value = -6350.3184
$ws.Range("K61").Value = 6562.3184
$ws.Range("N70").Value = -43140
$ws.Range("J70").Value = 42600
$ws.Range("L70").Value = 42600
$ws.Range("H70").Value = 42600
$ws.Range("N73").Value = -44472
$ws.Range("H73").Value = 42600
$ws.Range("L73").Value = 42600
$ws.Range("J73").Value = 42600
$ws.Range("M102").Value = -7407345
$ws.Range("K102").Value = 7408967
$ws.Range("I102").Value = 7408967
$ws.Range("L102").Value = 1550
$ws.Range("H102").Value = 6174397.5
$ws.Range("J102").Value = 1550
$ws.Range("N102").Value = -4794
$ws.Range("K136").Value = 19686.9552
$ws.Range("I136").Value = 6562.3184
$ws.Range("J136").Value = 3987.6843
$ws.Range("H136").Value = 5369.1953
$ws.Range("M136").Value = -17136.9552
$ws.Range("N136").Value = -17063.0529
$ws.Range("L136").Value = 11963.0529

# --- Sheet: BSM (23 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K20").Value = 2052.125
$ws.Range("I20").Value = 2052.125
$ws.Range("M20").Value = -1805.125
$ws.Range("H20").Value = 13487.954
$ws.Range("M46").ClearContents()
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("K99").Value = 142858030
$ws.Range("L99").Value = 1877.9
$ws.Range("I99").Value = 142858030
$ws.Range("J99").Value = 1877.9
$ws.Range("M99").Value = -142856532
$ws.Range("N99").Value = -4873.9
$ws.Range("H99").Value = 58824996
$ws.Range("N103").ClearContents()
$ws.Range("J103").Value = 0
$ws.Range("H103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("M107").Value = 950.3077
$ws.Range("I107").Value = 969.6923
$ws.Range("K107").Value = 969.6923
$ws.Range("H107").Value = 1152.7894

# --- Sheet: CRP (25 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J4").Value = 90002
$ws.Range("H4").Value = 90002
$ws.Range("N4").Value = -90226
$ws.Range("L4").Value = 90002
$ws.Range("H105").Value = 2174.9
$ws.Range("I105").Value = 2372.4167
$ws.Range("M105").Value = -625.4167000000002
$ws.Range("L105").Value = 1878.625
$ws.Range("K105").Value = 2372.4167
$ws.Range("N105").Value = -5372.625
$ws.Range("J105").Value = 1878.625
$ws.Range("N122").Value = -18083.3329
$ws.Range("J122").Value = 4394.4443
$ws.Range("K122").Value = 83333328
$ws.Range("H122").Value = 2781732.5
$ws.Range("I122").Value = 27777776
$ws.Range("L122").Value = 13183.3329
$ws.Range("M122").Value = -83330878
$ws.Range("H134").Value = 3236.8857
$ws.Range("N134").Value = -14386.2
$ws.Range("J134").Value = 3105.4
$ws.Range("M134").Value = -7241.400000000001
$ws.Range("K134").Value = 9776.400000000001
$ws.Range("L134").Value = 9316.200000000001
$ws.Range("I134").Value = 3258.8

# --- Sheet: CUL (14 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1111600
$ws.Range("M113").Value = -3408383.3
$ws.Range("N113").Value = -3005824.7
$ws.Range("L113").Value = 3001484.7
$ws.Range("I113").Value = 1136851.1
$ws.Range("K113").Value = 3410553.3
$ws.Range("J113").Value = 1000494.9
$ws.Range("J132").Value = 3500
$ws.Range("I132").Value = 2880
$ws.Range("N132").Value = -36560
$ws.Range("H132").Value = 3057.1428
$ws.Range("K132").Value = 25920
$ws.Range("L132").Value = 31500
$ws.Range("M132").Value = -23390

# --- Sheet: GSM (35 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N7").ClearContents()
$ws.Range("H7").Value = 25000000
$ws.Range("M7").Value = -24999888
$ws.Range("L7").Value = 0
$ws.Range("K7").Value = 25000000
$ws.Range("J7").Value = 0
$ws.Range("I7").Value = 25000000
$ws.Range("N8").ClearContents()
$ws.Range("I8").Value = 25000000
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("H8").Value = 25000000
$ws.Range("K8").Value = 25000000
$ws.Range("M8").Value = -24999861
$ws.Range("K80").Value = 2600
$ws.Range("I80").Value = 2600
$ws.Range("L80").Value = 3000
$ws.Range("N80").Value = -4996
$ws.Range("J80").Value = 3000
$ws.Range("M80").Value = -1602
$ws.Range("H80").Value = 2800
$ws.Range("M83").Value = -8008
$ws.Range("I83").Value = 2600
$ws.Range("L83").Value = 15000
$ws.Range("N83").Value = -24984
$ws.Range("H83").Value = 2800
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 13000
$ws.Range("J132").Value = 2701.625
$ws.Range("I132").Value = 5346.7
$ws.Range("N132").Value = -13164.875
$ws.Range("H132").Value = 3718.9614
$ws.Range("K132").Value = 16040.1
$ws.Range("L132").Value = 8104.875
$ws.Range("M132").Value = -13510.1

# --- Sheet: LTW (25 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N34").ClearContents()
$ws.Range("L34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L61").Value = 1255
$ws.Range("J61").Value = 1255
$ws.Range("N61").Value = -1659
$ws.Range("I61").Value = 1418.4445
$ws.Range("H61").Value = 1346.9375
$ws.Range("M61").Value = -1216.4445
$ws.Range("K61").Value = 1418.4445
$ws.Range("L93").Value = 900
$ws.Range("J93").Value = 900
$ws.Range("M93").Value = 648
$ws.Range("H93").Value = 787.5
$ws.Range("N93").Value = -3396
$ws.Range("I93").Value = 600
$ws.Range("K93").Value = 600
$ws.Range("H113").Value = 1346.9375
$ws.Range("M113").Value = 751.5554999999999
$ws.Range("N113").Value = -5595
$ws.Range("L113").Value = 1255
$ws.Range("I113").Value = 1418.4445
$ws.Range("K113").Value = 1418.4445
$ws.Range("J113").Value = 1255

# --- Sheet: WVR (28 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 943.8
$ws.Range("M113").Value = 128.2501
$ws.Range("N113").Value = -10329.9998
$ws.Range("L113").Value = 5989.9998
$ws.Range("I113").Value = 680.5833
$ws.Range("K113").Value = 2041.7499
$ws.Range("J113").Value = 1996.6666
$ws.Range("N122").Value = -9901.4998
$ws.Range("J122").Value = 1667.1666
$ws.Range("K122").Value = 3300.75
$ws.Range("H122").Value = 1440.4
$ws.Range("I122").Value = 1100.25
$ws.Range("L122").Value = 5001.4998
$ws.Range("M122").Value = -850.75
$ws.Range("J132").Value = 3483
$ws.Range("I132").Value = 3333.3333
$ws.Range("N132").Value = -15509
$ws.Range("H132").Value = 3433.111
$ws.Range("K132").Value = 9999.999899999999
$ws.Range("L132").Value = 10449
$ws.Range("M132").Value = -7469.999899999999
$ws.Range("K136").Value = 2082.0909
$ws.Range("I136").Value = 694.0303
$ws.Range("J136").Value = 1834.091
$ws.Range("H136").Value = 1150.0546
$ws.Range("M136").Value = 467.9090999999999
$ws.Range("N136").Value = -10602.273
$ws.Range("L136").Value = 5502.272999999999
